$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# ALC
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H106").Value = 2121.2856
$ws.Range("I106").Value = 799.75
$ws.Range("K106").Value = 799.75
$ws.Range("M106").Value = -168.75

$ws.Range("H132").Value = 93857.41
$ws.Range("I132").Value = 69791.266
$ws.Range("J132").Value = 145427.72
$ws.Range("K132").Value = 209373.798
$ws.Range("L132").Value = 436283.16
$ws.Range("M132").Value = -206843.798
$ws.Range("N132").Value = -441343.16

$ws.Range("H137").Value = 1246362.6
$ws.Range("I137").Value = 24760.68
$ws.Range("J137").Value = 5063868.5
$ws.Range("K137").Value = 74282.04000000001
$ws.Range("L137").Value = 15191605.5
$ws.Range("M137").Value = -71732.04000000001
$ws.Range("N137").Value = -15196705.5

$ws.Range("H138").Value = 4242.092
$ws.Range("I138").Value = 2384.8
$ws.Range("J138").Value = 4483.299
$ws.Range("K138").Value = 7154.400000000001
$ws.Range("L138").Value = 13449.897
$ws.Range("M138").Value = -2014.400000000001
$ws.Range("N138").Value = -23729.897

$ws.Range("H140").Value = 65048
$ws.Range("J140").Value = 65048
$ws.Range("L140").Value = 65048
$ws.Range("N140").Value = -75408

# ---------------------------------------------------------------------------
# ARM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H32").Value = 3600.36
$ws.Range("I32").Value = 2672.5476
$ws.Range("J32").Value = 8471.375
$ws.Range("K32").Value = 2672.5476
$ws.Range("L32").Value = 8471.375
$ws.Range("M32").Value = -2385.5476
$ws.Range("N32").Value = -9045.375

$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()

$ws.Range("H45").Value = 24848.715
$ws.Range("I45").Value = 24177
$ws.Range("K45").Value = 24177
$ws.Range("M45").Value = -23800

$ws.Range("H132").Value = 3898.9443
$ws.Range("I132").Value = 3646.9656
$ws.Range("K132").Value = 10940.8968
$ws.Range("M132").Value = -8410.8968

# ---------------------------------------------------------------------------
# BSM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H10").Value = 330
$ws.Range("I10").Value = 330
$ws.Range("K10").Value = 330
$ws.Range("M10").Value = -190

$ws.Range("H12").Value = 671
$ws.Range("I12").Value = 205
$ws.Range("J12").Value = 787.5
$ws.Range("K12").Value = 205
$ws.Range("L12").Value = 787.5
$ws.Range("M12").Value = -37
$ws.Range("N12").Value = -1123.5

$ws.Range("H81").Value = 30234.285
$ws.Range("J81").Value = 30234.285
$ws.Range("L81").Value = 30234.285
$ws.Range("N81").Value = -32356.285

$ws.Range("H84").Value = 30234.285
$ws.Range("J84").Value = 30234.285
$ws.Range("L84").Value = 90702.855
$ws.Range("N84").Value = -101310.855

$ws.Range("H133").Value = 65000
$ws.Range("J133").Value = 65000
$ws.Range("L133").Value = 65000
$ws.Range("N133").Value = -75120

$ws.Range("H134").Value = 1769.1072
$ws.Range("I134").Value = 1524.4231
$ws.Range("K134").Value = 4573.2693
$ws.Range("M134").Value = -2038.2693

$ws.Range("H135").Value = 51325
$ws.Range("J135").Value = 51325
$ws.Range("L135").Value = 51325
$ws.Range("N135").Value = -61465

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws.Range("H138").Value = 99878.5
$ws.Range("J138").Value = 99878.5
$ws.Range("L138").Value = 99878.5
$ws.Range("N138").Value = -110158.5

$ws.Range("H139").Value = 99989
$ws.Range("J139").Value = 99989
$ws.Range("L139").Value = 99989
$ws.Range("N139").Value = -110269

$ws.Range("H140").Value = 77737
$ws.Range("J140").Value = 77737
$ws.Range("L140").Value = 77737
$ws.Range("N140").Value = -88097

# ---------------------------------------------------------------------------
# CRP
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H118").Value = 129000
$ws.Range("J118").Value = 129000
$ws.Range("L118").Value = 129000
$ws.Range("N118").Value = -132314

$ws.Range("H140").Value = 79957
$ws.Range("J140").Value = 79957
$ws.Range("L140").Value = 79957
$ws.Range("N140").Value = -90317

# ---------------------------------------------------------------------------
# CUL
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H68").Value = 1472422.9
$ws.Range("I68").Value = 1696
$ws.Range("J68").Value = 1725996.5
$ws.Range("K68").Value = 5088
$ws.Range("L68").Value = 5177989.5
$ws.Range("M68").Value = -4277
$ws.Range("N68").Value = -5179611.5

$ws.Range("H71").Value = 1472422.9
$ws.Range("I71").Value = 1696
$ws.Range("J71").Value = 1725996.5
$ws.Range("K71").Value = 15264
$ws.Range("L71").Value = 15533968.5
$ws.Range("M71").Value = -11208
$ws.Range("N71").Value = -15542080.5

$ws.Range("H118").Value = 2921.125
$ws.Range("I118").Value = 2476.7144
$ws.Range("K118").Value = 7430.1432
$ws.Range("M118").Value = -6187.1432

$ws.Range("H119").Value = 1002.25
$ws.Range("I119").Value = 1002.25
$ws.Range("K119").Value = 3006.75
$ws.Range("M119").Value = 1831.25

$ws.Range("H120").Value = 22198.75
$ws.Range("I120").Value = 17518
$ws.Range("K120").Value = 52554
$ws.Range("M120").Value = -47716

$ws.Range("H129").Value = 4126835
$ws.Range("I129").Value = 7616063.5
$ws.Range("K129").Value = 22848190.5
$ws.Range("M129").Value = -22843190.5

$ws.Range("H141").Value = 2745.7144
$ws.Range("I141").Value = 2364.5
$ws.Range("J141").Value = 5033
$ws.Range("K141").Value = 7093.5
$ws.Range("L141").Value = 15099
$ws.Range("M141").Value = -1913.5
$ws.Range("N141").Value = -25459

# ---------------------------------------------------------------------------
# GSM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H39").Value = 31261
$ws.Range("J39").Value = 31261
$ws.Range("L39").Value = 31261
$ws.Range("N39").Value = -32325

$ws.Range("H126").Value = 22612.941
$ws.Range("J126").Value = 5128.4287
$ws.Range("L126").Value = 15385.2861
$ws.Range("N126").Value = -20325.2861

$ws.Range("H132").Value = 24058.695
$ws.Range("I132").Value = 24058.695
$ws.Range("K132").Value = 72176.08499999999
$ws.Range("M132").Value = -69646.08499999999

# ---------------------------------------------------------------------------
# LTW
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H61").Value = 3892.5676
$ws.Range("I61").Value = 3728
$ws.Range("J61").Value = 4235.4165
$ws.Range("K61").Value = 3728
$ws.Range("L61").Value = 4235.4165
$ws.Range("M61").Value = -3526
$ws.Range("N61").Value = -4639.4165

$ws.Range("H113").Value = 3892.5676
$ws.Range("I113").Value = 3728
$ws.Range("J113").Value = 4235.4165
$ws.Range("K113").Value = 3728
$ws.Range("L113").Value = 4235.4165
$ws.Range("M113").Value = -1558
$ws.Range("N113").Value = -8575.416499999999

$ws.Range("H132").Value = 4909.4287
$ws.Range("I132").Value = 4074
$ws.Range("J132").Value = 6023.3335
$ws.Range("K132").Value = 12222
$ws.Range("L132").Value = 18070.0005
$ws.Range("M132").Value = -9692
$ws.Range("N132").Value = -23130.0005

$ws.Range("H141").Value = 90650
$ws.Range("J141").Value = 90650
$ws.Range("L141").Value = 90650
$ws.Range("N141").Value = -101010

# ---------------------------------------------------------------------------
# WVR
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H122").Value = 4436.467
$ws.Range("I122").Value = 2699.739
$ws.Range("J122").Value = 10142.857
$ws.Range("K122").Value = 8099.217000000001
$ws.Range("L122").Value = 30428.571
$ws.Range("M122").Value = -5649.217000000001
$ws.Range("N122").Value = -35328.571

$ws.Range("H132").Value = 2953.5
$ws.Range("I132").Value = 3058
$ws.Range("J132").Value = 550
$ws.Range("K132").Value = 9174
$ws.Range("L132").Value = 1650
$ws.Range("M132").Value = -6644
$ws.Range("N132").Value = -6710

$ws.Range("H138").Value = 84774.5
$ws.Range("J138").Value = 79550
$ws.Range("L138").Value = 79550
$ws.Range("N138").Value = -89830

$ws.Range("H139").Value = 79626.664
$ws.Range("J139").Value = 79626.664
$ws.Range("L139").Value = 79626.664
$ws.Range("N139").Value = -89906.664
